# "Generate Report for Handoff" — refresh the localization-status report with
# a new handoff round: new source-file GUIDs, new handoff timestamps, status
# flips from "Handed back" to "Ready for handoff", and the (now-empty)
# "Latest Target File" / "Latest Handback File" columns are cleared because
# nothing has been handed back yet for this round.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-45-12 10:45:43"

$ov.Range("A3").Value = "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-45-12 10:45:43"

# Recreate the hyperlinks on column A so their display text matches the new
# file names (deleting via a range wipes every hyperlink on the sheet in
# this runtime, so we delete once and re-add everything that should stay).
$ov.Range("A2").Hyperlinks.Delete()

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/a1e216f4-665e-4646-a68d-b253891e61e3.md", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/ffff607a3117-83cd-468a-8a52-3bd8918277dd.md", $null, $null, "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md") | Out-Null

Write-Host "Overview sheet updated"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-12 10:45:39"
$zh.Range("H2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-12 10:45:39"
$zh.Range("H3").Value = "0001-01-01 00:00:00"

# Latest Target File / Latest Handback File: nothing has been handed back
# yet for this round, so these columns go blank.
$zh.Range("F2:G3").Clear()

# Rebuild hyperlinks: deleting from a range wipes the whole sheet's
# hyperlink collection in this runtime, so delete once then re-add every
# link that should survive (F/G links are intentionally NOT re-added).
$zh.Range("A2").Hyperlinks.Delete()

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/a1e216f4-665e-4646-a68d-b253891e61e3.md", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/a1e216f4-665e-4646-a68d-b253891e61e3.md", $null, $null, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4cd6174898dcc6d1eafe54a0cae57f36067e22b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/ffff607a3117-83cd-468a-8a52-3bd8918277dd.md", $null, $null, "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/ffff607a3117-83cd-468a-8a52-3bd8918277dd.md", $null, $null, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4cd6174898dcc6d1eafe54a0cae57f36067e22b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.zh-cn.xlf") | Out-Null

Write-Host "zh-cn sheet updated"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf"
$de.Range("E2").Value = "2016-03-12 10:45:43"
$de.Range("H2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf"
$de.Range("E3").Value = "2016-03-12 10:45:43"
$de.Range("H3").Value = "0001-01-01 00:00:00"

# Latest Target File / Latest Handback File: nothing has been handed back
# yet for this round, so these columns go blank.
$de.Range("F2:G3").Clear()

# Rebuild hyperlinks the same way as the zh-cn sheet.
$de.Range("A2").Hyperlinks.Delete()

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/a1e216f4-665e-4646-a68d-b253891e61e3.md", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/a1e216f4-665e-4646-a68d-b253891e61e3.md", $null, $null, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b861051ec2ffc55e74a98923d07f90eea17aa80a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/ffff607a3117-83cd-468a-8a52-3bd8918277dd.md", $null, $null, "ffff607a3117-83cd-468a-8a52-3bd8918277dd.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5658a74b205255cfdb195dfd3fc8335f0ad5fc74/e2e/ffff607a3117-83cd-468a-8a52-3bd8918277dd.md", $null, $null, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b861051ec2ffc55e74a98923d07f90eea17aa80a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf", $null, $null, "a1e216f4-665e-4646-a68d-b253891e61e3.26121daf7d17a45b4bebbc4c245f3629c07902c0.de-de.xlf") | Out-Null

Write-Host "de-de sheet updated"
